$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.71"
$ws.Range("E2").Value = "'-4.02%"
$ws.Range("D3").Value = "'40.10"
$ws.Range("E3").Value = "'-5.97%"
$ws.Range("D4").Value = "'5.089"
$ws.Range("E4").Value = "'-2.03%"
$ws.Range("D5").Value = "'0.07686"
$ws.Range("E5").Value = "'-6.26%"
$ws.Range("D6").Value = "'4.243"
$ws.Range("E6").Value = "'-2.52%"
$ws.Range("D7").Value = "'1.622"
$ws.Range("E7").Value = "'-8.41%"
$ws.Range("D8").Value = "'0.9172"
$ws.Range("E8").Value = "'-3.56%"
$ws.Range("D9").Value = "'0.1038"
$ws.Range("E9").Value = "'-7.00%"
$ws.Range("D10").Value = "'0.1778"
$ws.Range("E10").Value = "'-5.69%"
$ws.Range("D11").Value = "'0.09298"
$ws.Range("E11").Value = "'-1.33%"
$ws.Range("D12").Value = "'0.04438"
$ws.Range("E12").Value = "'-5.17%"
$ws.Range("D13").Value = "'0.1055"
$ws.Range("E13").Value = "'-0.20%"
$ws.Range("D14").Value = "'0.001262"
$ws.Range("E14").Value = "'-3.17%"
$ws.Range("D15").Value = "'0.005847"
$ws.Range("E15").Value = "'1.68%"
$ws.Range("E16").Value = "'2,409.05%"
$ws.Range("D17").Value = "'3.362"
$ws.Range("E17").Value = "'0.22%"
$ws.Range("E18").Value = "'-4.37%"
$ws.Range("D19").Value = "'0.3313"
$ws.Range("E19").Value = "'-1.57%"
$ws.Range("D20").Value = "'6.907"
$ws.Range("E20").Value = "'-7.29%"
$ws.Range("E21").Value = "'-2.05%"
$ws.Range("D22").Value = "'0.2704"
$ws.Range("E22").Value = "'5.87%"
$ws.Range("D23").Value = "'0.04152"
$ws.Range("E23").Value = "'-0.78%"
$ws.Range("D24").Value = "'0.001202"
$ws.Range("E24").Value = "'-3.55%"
$ws.Range("D25").Value = "'0.004113"
$ws.Range("E25").Value = "'-3.82%"
$ws.Range("E26").Value = "'6.17%"
$ws.Range("D38").Value = "'0.02488"
$ws.Range("E38").Value = "'-6.47%"
$ws.Range("D39").Value = "'0.05197"
$ws.Range("E39").Value = "'-8.07%"
$ws.Range("D40").Value = "'0.007930"
$ws.Range("E40").Value = "'-2.44%"
$ws.Range("D41").Value = "'0.1317"
$ws.Range("E41").Value = "'-6.12%"
$ws.Range("D42").Value = "'0.007075"
$ws.Range("E42").Value = "'7.71%"
$ws.Range("E43").Value = "'-4.33%"
$ws.Range("D44").Value = "'0.007416"
$ws.Range("E44").Value = "'-3.74%"
$ws.Range("D45").Value = "'0.3069"
$ws.Range("E45").Value = "'-11.91%"
$ws.Range("D46").Value = "'0.00006406"
$ws.Range("E46").Value = "'-5.36%"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.27%"
$ws.Range("D48").Value = "'0.002999"
$ws.Range("E48").Value = "'-27.03%"
$ws.Range("D49").Value = "'0.004534"
$ws.Range("E49").Value = "'35.36%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.27%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.27%"
